$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Genius / Him & I mismatch
$ws.Range("A6").Value = "Genius"
$ws.Range("B6").Value = "N"
$ws.Range("D6").Value = "http://genius.com/G-Eazy-Him-&-I-lyrics"
$ws.Range("E6").Value = "https://genius.com/G-eazy-him-and-i-lyrics"
$ws.Range("F6").Value = '"and" instead of "&"'
$ws.Range("C6").Value = "Him & I"

# Row 7: Genius / Flasker på flasker mismatch
$ws.Range("A7").Value = "Genius"
$ws.Range("B7").Value = "N"
$ws.Range("D7").Value = "http://genius.com/Sushi-x-Kobe-Flasker-på-flasker-lyrics"
$ws.Range("E7").Value = "https://genius.com/Sushi-x-kobe-flasker-pa-flasker-lyrics"
$ws.Range("F7").Value = '"a" instead of "å"'
$ws.Range("C7").Value = "Flasker på flasker"

# Adjust column E width to fit new (wider) content
# (Target stored width 53.42578125 comes from Excel's real font-metrics based
#  bestFit; this headless engine quantizes ColumnWidth to coarser steps, so we
#  pick the input that lands on the closest achievable stored width.)
$ws.Columns.Item(5).ColumnWidth = 52.67

# Update selection to reflect where the user ended up (next empty row, col A)
$ws.Range("A8").Select()
